# ANS-006-150 - Verification du QA (#385)
# Regenerated StructureDefinition export: bump the generation timestamp and
# add a new "Context" row documenting the QuestionnaireResponse extension
# context on the Metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the generation Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T16:42:45+00:00"

# Append a new Context row (element:QuestionnaireResponse) right after the
# existing element:Encounter context row (row 21). Copy the row's
# formatting first so the new cells keep the same style as their peers,
# then overwrite the values.
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Context"
$ws.Range("B22").Value = "element:QuestionnaireResponse"
